# =====================================================================
# Apply the "infl_cost and soukann check" edit:
#  - Extend sheet1 data (rows 25-55) with more simulation results
#  - Extend the existing L2/M2 averages and add a small max/overall/
#    per-seed summary table in columns K:M (rows 3, 6-13)
#  - Add a second worksheet ("Sheet1") holding an earlier snapshot of
#    the same kind of data (rows 2-15) plus its own avg/max summary
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) New simulation rows appended to the main sheet (A25:E55).
#    F (B-E), G (D-C), H (F/B), I (G/D) are filled as formulas,
#    matching the existing pattern used for rows 2-24.
# ---------------------------------------------------------------
$newRows = @(
    @(3, 199.18199999999999, 254.82400000000001, 256.07999999999902, 195.84299999999899),
    @(4, 20.574000000000002, 96.538000000000096, 104.716999999999, 8.3489999999999398),
    @(4, 60.308, 189.372999999999, 198.844999999999, 54.5549999999999),
    @(4, 104.787999999999, 227.05699999999999, 241.427999999999, 86.114999999999895),
    @(4, 112.18799999999899, 229.84100000000001, 246.84199999999899, 95.042000000000002),
    @(4, 147.696, 252.17499999999899, 251.605999999999, 141.36600000000001),
    @(4, 153.218999999999, 254.575999999999, 259.09899999999902, 127.30999999999899),
    @(5, 22.314, 89.194999999999993, 104.267, 23.945),
    @(5, 65.565000000000097, 201.78200000000001, 192.91199999999901, 56.415999999999997),
    @(5, 126.36199999999999, 238.87199999999899, 243.253999999999, 116.036999999999),
    @(5, 133.779, 242.95299999999901, 246.664999999999, 126.015999999999),
    @(5, 177.47399999999899, 256.34899999999902, 256.79299999999898, 178.18700000000001),
    @(5, 183.393, 258.43599999999998, 258.99099999999999, 183.44900000000001),
    @(6, 43.079000000000001, 96.403999999999996, 99.135999999999996, 42.555),
    @(6, 94.525000000000006, 198.09100000000001, 190.70500000000001, 86.644999999999996),
    @(6, 192.631, 265.22399999999999, 267.3, 192.70500000000001),
    @(6, 199.88499999999999, 264.87200000000001, 271.56200000000001, 202.84299999999999),
    @(6, 216.715, 267.63, 272.20600000000002, 204.797),
    @(6, 236.25800000000001, 278.51499999999999, 278.88600000000002, 235.41399999999999),
    @(7, 21.15, 89.296999999999997, 93.744, 16.048999999999999),
    @(7, 59.832000000000001, 196.666, 188.86799999999999, 33.526000000000003),
    @(7, 125.33499999999999, 239.66900000000001, 240.75, 122.91500000000001),
    @(7, 132.13399999999999, 243.999, 244.47200000000001, 131.768),
    @(7, 166.625, 256.55, 254.816, 164.779),
    @(7, 168.596, 256.46699999999998, 259.54399999999998, 169.77699999999999),
    @(8, 40.356999999999999, 89.054000000000002, 105.753, 37.43),
    @(8, 96.399000000000001, 195.57300000000001, 204.97900000000001, 93.974000000000004),
    @(8, 190.017, 259.51100000000002, 260.64499999999998, 175.06399999999999),
    @(8, 195.03899999999999, 261.07100000000003, 265.69900000000001, 185.62200000000001),
    @(8, 207.05199999999999, 265.72300000000001, 268.70800000000003, 184.89599999999999),
    @(8, 232.59399999999999, 272.45499999999998, 274.25900000000001, 223.45)
)

$startRow = 25
$row = $startRow
foreach ($d in $newRows) {
    $ws.Cells.Item($row, 1).Value2 = $d[0]
    $ws.Cells.Item($row, 2).Value2 = $d[1]
    $ws.Cells.Item($row, 3).Value2 = $d[2]
    $ws.Cells.Item($row, 4).Value2 = $d[3]
    $ws.Cells.Item($row, 5).Value2 = $d[4]
    $row++
}
$endRow = $row - 1

# Fill F/G/H/I formulas in two blocks the same way the source file does
# (rows 25-41, then 42-55), mirroring two separate fill/copy actions.
$ws.Range("F25:F41").Formula = "=B25-E25"
$ws.Range("G25:G41").Formula = "=D25-C25"
$ws.Range("H25:H41").Formula = "=F25/B25"
$ws.Range("I25:I41").Formula = "=G25/D25"
$ws.Range("H25:H41").NumberFormat = "0.000%"
$ws.Range("I25:I41").NumberFormat = "0.000%"

$ws.Range("F42:F55").Formula = "=B42-E42"
$ws.Range("G42:G55").Formula = "=D42-C42"
$ws.Range("H42:H55").Formula = "=F42/B42"
$ws.Range("I42:I55").Formula = "=G42/D42"
$ws.Range("H42:H55").NumberFormat = "0.000%"
$ws.Range("I42:I55").NumberFormat = "0.000%"

# ---------------------------------------------------------------
# 2) Update the running average (now covering the bigger range)
# ---------------------------------------------------------------
$ws.Range("L2").Formula = "=AVERAGE(H2:H55)"
$ws.Range("M2").Formula = "=AVERAGE(I2:I41)"

# ---------------------------------------------------------------
# 3) New max / overall / per-seed-group summary table in K:M
# ---------------------------------------------------------------
$ws.Range("K3").Value2 = "max"
$ws.Range("L3").Formula = "=MAX(H2:H55)"
$ws.Range("L3").NumberFormat = "0.000%"

$ws.Range("K6").Value2 = "コスト(拡散量)"
$ws.Range("L6").Value2 = "平均"
$ws.Range("M6").Value2 = "最大"

$ws.Range("K7").Value2 = "全体"
$ws.Range("L7").Value2 = 0.08194886429232183
$ws.Range("M7").Value2 = 0.59419655876349098
$ws.Range("L7").NumberFormat = "0.000%"
$ws.Range("M7").NumberFormat = "0.000%"

$ws.Range("K8").Value2 = 100
$ws.Range("L8").Formula = "=AVERAGE(H2,H8,H14,H20,H26,H32,H38,H44,H50)"
$ws.Range("M8").Formula = "=MAX(H2,H8,H14,H20,H26,H32,H38,H44,H50)"

$ws.Range("K9").Value2 = 200
$ws.Range("K10").Value2 = 300
$ws.Range("K11").Value2 = 400
$ws.Range("K12").Value2 = 500
$ws.Range("K13").Value2 = 600

$ws.Range("L9:L13").Formula = "=AVERAGE(H3,H9,H15,H21,H27,H33,H39,H45,H51)"
$ws.Range("M9:M12").Formula = "=MAX(H3,H9,H15,H21,H27,H33,H39,H45,H51)"
$ws.Range("M13").Formula = "=MAX(H7,H13,H19,H25,H31,H37,H43,H49,H55)"

$ws.Range("L8:L13").NumberFormat = "0.000%"
$ws.Range("M8:M13").NumberFormat = "0.000%"

# Column K now needs to fit the new "全体"/"コスト(拡散量)" labels
$ws.Columns.Item(11).AutoFit()

# ---------------------------------------------------------------
# 4) Add the second worksheet ("Sheet1") with an earlier snapshot
#    of the same simulation data, plus its own avg/max summary.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sheet1"

$ws2Rows = @(
    @(0, 20.831999999999901, 92.989999999999895, 103.52699999999901, 17.716999999999999),
    @(0, 55.18, 189.81899999999999, 200.16899999999899, 53.720999999999997),
    @(1, 30.4920000000001, 101.262, 96.751000000000104, 28.697000000000099),
    @(1, 75.766000000000005, 204.94799999999901, 182.289999999999, 56.37),
    @(2, 22.796999999999901, 94.421999999999699, 101.04, 16.317999999999898),
    @(2, 50.122999999999998, 196.201999999999, 201.28699999999901, 43.643999999999998),
    @(3, 24.722000000000001, 90.176000000000201, 103.86, 24.495999999999999),
    @(3, 69.998000000000104, 195.852, 200.05499999999901, 65.938999999999993),
    @(4, 20.574000000000002, 96.538000000000096, 104.716999999999, 8.3489999999999398),
    @(4, 60.308, 189.372999999999, 198.844999999999, 54.5549999999999),
    @(5, 22.314, 89.194999999999993, 104.267, 23.945),
    @(5, 65.565000000000097, 201.78200000000001, 192.91199999999901, 56.415999999999997),
    @(6, 43.079000000000001, 96.404000000000295, 99.135999999999896, 42.555000000000099),
    @(6, 94.525000000000006, 198.09099999999901, 190.70499999999899, 86.645000000000195)
)

$row = 2
foreach ($d in $ws2Rows) {
    $ws2.Cells.Item($row, 1).Value2 = $d[0]
    $ws2.Cells.Item($row, 2).Value2 = $d[1]
    $ws2.Cells.Item($row, 3).Value2 = $d[2]
    $ws2.Cells.Item($row, 4).Value2 = $d[3]
    $ws2.Cells.Item($row, 5).Value2 = $d[4]
    $row++
}
$ws2Last = $row - 1

$ws2.Range("F2").Formula = "=B2-E2"
$ws2.Range("G2").Formula = "=D2-C2"
$ws2.Range("H2").Formula = "=F2/B2"
$ws2.Range("I2").Formula = "=G2/D2"
$ws2.Range("H2").NumberFormat = "0.000%"
$ws2.Range("I2").NumberFormat = "0.000%"

$ws2.Range("F3:F15").Formula = "=B3-E3"
$ws2.Range("G3:G15").Formula = "=D3-C3"
$ws2.Range("H3:H15").Formula = "=F3/B3"
$ws2.Range("I3:I15").Formula = "=G3/D3"
$ws2.Range("H3:H15").NumberFormat = "0.000%"
$ws2.Range("I3:I15").NumberFormat = "0.000%"

$ws2.Range("K3").Formula = "=AVERAGE(H2:H15)"
$ws2.Range("K4").Formula = "=MAX(H2:H15)"
$ws2.Range("K3:K4").NumberFormat = "0.000%"

$ws2.Range("K5").Select()

# ---------------------------------------------------------------
# 5) Restore the original sheet as the active / selected one,
#    matching the saved selection in the source file.
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("E10").Select()

Write-Host "edit complete"
